{"js": "const replacements = [\n  { find: \"2025-10-20 Monday\", repl: \"2025-10-21 Tuesday\" },\n  { find: \"996\u00f77=142, 2\", repl: \"314\u00f77=44, 6\" },\n  { find: \"794\u00f78=99, 2\", repl: \"454\u00f79=50, 4\" },\n  { find: \"524\u00f78=65, 4\", repl: \"509\u00f79=56, 5\" },\n  { find: \"806\u00f73=268, 2\", repl: \"319\u00f75=63, 4\" },\n  { find: \"330\u00f79=36, 6\", repl: \"919\u00f73=306, 1\" },\n  { find: \"645\u00f79=71, 6\", repl: \"563\u00f77=80, 3\" },\n  { find: \"193\u00f79=21, 4\", repl: \"821\u00f79=91, 2\" },\n  { find: \"540\u00f79=60, 0\", repl: \"271\u00f74=67, 3\" },\n  { find: \"672\u00f79=74, 6\", repl: \"320\u00f76=53, 2\" },\n  { find: \"173\u00f74=43, 1\", repl: \"360\u00f73=120, 0\" },\n  { find: \"958\u00f75=191, 3\", repl: \"456\u00f72=228, 0\" },\n  { find: \"338\u00f74=84, 2\", repl: \"816\u00f72=408, 0\" },\n  { find: \"360\u00f75=72, 0\", repl: \"716\u00f72=358, 0\" },\n  { find: \"200\u00f72=100, 0\", repl: \"709\u00f76=118, 1\" },\n  { find: \"384\u00f72=192, 0\", repl: \"927\u00f74=231, 3\" },\n  { find: \"552\u00f75=110, 2\", repl: \"157\u00f78=19, 5\" },\n  { find: \"122\u00f73=40, 2\", repl: \"770\u00f76=128, 2\" },\n  { find: \"111\u00f75=22, 1\", repl: \"696\u00f72=348, 0\" },\n  { find: \"801\u00f74=200, 1\", repl: \"419\u00f77=59, 6\" },\n  { find: \"404\u00f75=80, 4\", repl: \"474\u00f76=79, 0\" },\n  { find: \"187\u00f77=26, 5\", repl: \"453\u00f74=113, 1\" },\n  { find: \"489\u00f75=97, 4\", repl: \"683\u00f77=97, 4\" },\n  { find: \"461\u00f75=92, 1\", repl: \"497\u00f74=124, 1\" },\n  { find: \"694\u00f73=231, 1\", repl: \"385\u00f78=48, 1\" },\n  { find: \"538\u00f77=76, 6\", repl: \"664\u00f78=83, 0\" },\n];\n\nconst body = context.document.body;\nfor (const { find, repl } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(repl, Word.InsertLocation.replace);\n  }\n  if (results.items.length === 0) {\n    console.log(`warning: no match found for \"${find}\"`);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"2025-10-20 Monday\"; Replace = \"2025-10-21 Tuesday\" },\n    @{ Find = \"996\u00f77=142, 2\"; Replace = \"314\u00f77=44, 6\" },\n    @{ Find = \"794\u00f78=99, 2\"; Replace = \"454\u00f79=50, 4\" },\n    @{ Find = \"524\u00f78=65, 4\"; Replace = \"509\u00f79=56, 5\" },\n    @{ Find = \"806\u00f73=268, 2\"; Replace = \"319\u00f75=63, 4\" },\n    @{ Find = \"330\u00f79=36, 6\"; Replace = \"919\u00f73=306, 1\" },\n    @{ Find = \"645\u00f79=71, 6\"; Replace = \"563\u00f77=80, 3\" },\n    @{ Find = \"193\u00f79=21, 4\"; Replace = \"821\u00f79=91, 2\" },\n    @{ Find = \"540\u00f79=60, 0\"; Replace = \"271\u00f74=67, 3\" },\n    @{ Find = \"672\u00f79=74, 6\"; Replace = \"320\u00f76=53, 2\" },\n    @{ Find = \"173\u00f74=43, 1\"; Replace = \"360\u00f73=120, 0\" },\n    @{ Find = \"958\u00f75=191, 3\"; Replace = \"456\u00f72=228, 0\" },\n    @{ Find = \"338\u00f74=84, 2\"; Replace = \"816\u00f72=408, 0\" },\n    @{ Find = \"360\u00f75=72, 0\"; Replace = \"716\u00f72=358, 0\" },\n    @{ Find = \"200\u00f72=100, 0\"; Replace = \"709\u00f76=118, 1\" },\n    @{ Find = \"384\u00f72=192, 0\"; Replace = \"927\u00f74=231, 3\" },\n    @{ Find = \"552\u00f75=110, 2\"; Replace = \"157\u00f78=19, 5\" },\n    @{ Find = \"122\u00f73=40, 2\"; Replace = \"770\u00f76=128, 2\" },\n    @{ Find = \"111\u00f75=22, 1\"; Replace = \"696\u00f72=348, 0\" },\n    @{ Find = \"801\u00f74=200, 1\"; Replace = \"419\u00f77=59, 6\" },\n    @{ Find = \"404\u00f75=80, 4\"; Replace = \"474\u00f76=79, 0\" },\n    @{ Find = \"187\u00f77=26, 5\"; Replace = \"453\u00f74=113, 1\" },\n    @{ Find = \"489\u00f75=97, 4\"; Replace = \"683\u00f77=97, 4\" },\n    @{ Find = \"461\u00f75=92, 1\"; Replace = \"497\u00f74=124, 1\" },\n    @{ Find = \"694\u00f73=231, 1\"; Replace = \"385\u00f78=48, 1\" },\n    @{ Find = \"538\u00f77=76, 6\"; Replace = \"664\u00f78=83, 0\" },\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # Wrap:=1 (wdFindContinue), Replace:=2 (wdReplaceAll)\n    $find.Execute($r.Find, $true, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2) | Out-Null\n}"}
